# Auto update Excel log
# Appends new sensor log rows to the ALERTS sheet (rows 14-15) and the
# mmWave sheet (rows 59-62), matching the upstream log export.

$wb = $excel.ActiveWorkbook

# --- ALERTS sheet: two new FALL_DETECTED rows ---
$alerts = $wb.Worksheets.Item("ALERTS")

$alertsData = @(
    @("2026-02-01", "11:25:29", "11:00", "Living Room", "CRITICAL", "FALL_DETECTED"),
    @("2026-02-01", "11:25:32", "11:00", "Living Room", "CRITICAL", "FALL_DETECTED")
)

$r = 14
foreach ($row in $alertsData) {
    # Columns A-F hold plain text in this log (dates/times are literal
    # strings, not Excel date/time serials). Force text format on the date
    # cell so Excel doesn't auto-convert "2026-02-01" into a date serial,
    # then clear the formatting so the cell keeps the sheet's default
    # (unstyled) look, matching the rest of the log.
    $dateCell = $alerts.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row[0]
    $dateCell.ClearFormats()

    $alerts.Cells.Item($r, 2).Value = $row[1]
    $alerts.Cells.Item($r, 3).Value = $row[2]
    $alerts.Cells.Item($r, 4).Value = $row[3]
    $alerts.Cells.Item($r, 5).Value = $row[4]
    $alerts.Cells.Item($r, 6).Value = $row[5]

    $r = $r + 1
}

# --- mmWave sheet: four new PRESENCE_DETECTED rows ---
$mmwave = $wb.Worksheets.Item("mmWave")

$mmwaveData = @(
    @("2026-02-01", "11:25:32", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:25:42", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:25:53", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:26:03", "11:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$r = 59
foreach ($row in $mmwaveData) {
    $dateCell = $mmwave.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row[0]
    $dateCell.ClearFormats()

    $mmwave.Cells.Item($r, 2).Value = $row[1]
    $mmwave.Cells.Item($r, 3).Value = $row[2]
    $mmwave.Cells.Item($r, 4).Value = $row[3]
    $mmwave.Cells.Item($r, 5).Value = $row[4]
    $mmwave.Cells.Item($r, 6).Value = $row[5]

    $r = $r + 1
}
